# Insert a new data row at row 117 (shifting existing rows 117-135 down to 118-136)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 117..135 down by one row, creating a new blank row 117.
$ws.Rows("117:117").Insert()

# Populate the newly inserted row 117 with the new record.
$ws.Range("A117").Value = 7
$ws.Range("B117").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C117").Value = "Ñuble"
$ws.Range("D117").Value = 44491
$ws.Range("D117").NumberFormat = $ws.Range("D118").NumberFormat
$ws.Range("E117").Value = 16
$ws.Range("F117").Value = 100112006
$ws.Range("G117").Value = "Repollo"
$ws.Range("H117").Value = "Crespo record"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 600
$ws.Range("L117").Value = 700
$ws.Range("M117").Value = 650
$ws.Range("N117").Value = "$/unidad"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 650
$ws.Range("Q117").Value = 1
$ws.Range("R117").Value = "Hortaliza"
